# Append: 2025-10-14 01:43 JST
# A new listing ("Drupal..." ) was scraped and inserted as row 13, pushing
# the two previously-last listings down by one row. Every already-present
# row also gets its "取得日時" (fetched-at) timestamp refreshed to the new
# run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-14 01:43:52"

# 1) Insert a fresh blank row at 13; rows 13-14 (the two oldest listings)
#    shift down to 14-15, carrying their values/styles/formatting with them.
$ws.Rows.Item(13).Insert()

# 2) Refresh the "fetched at" timestamp on every row that already existed
#    (rows 2-12 stayed in place; rows 14-15 are the ones that just shifted
#    down from 13-14).
foreach ($r in 2..12) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
foreach ($r in 14..15) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# 3) Populate the brand-new row 13 with the newly scraped listing.
$ws.Cells.Item(13, 1).Value = $newTimestamp
$ws.Cells.Item(13, 2).Value = "Drupal関連プロジェクトの要件定義や基本設計ができる方(1人月、長期継続案件)"
$ws.Cells.Item(13, 3).Value = "システム開発"
$ws.Cells.Item(13, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(13, 5).Value = "期限情報なし"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5400683"
$ws.Cells.Item(13, 7).Value = 25

# 4) The plain row-insert leaves the worksheet's <hyperlinks> list pointing
#    at stale cell references (it shifts the underlying relationship
#    targets but not the recorded ref= cell), and the newly inserted F13
#    has no hyperlink entry at all yet. Rebuild the hyperlink list from
#    scratch, in row order, so every F-column URL cell (F2:F15) gets
#    exactly one correct hyperlink pointing at its own text.
$urls = @{
    2  = "https://www.lancers.jp/work/detail/5412417"
    3  = "https://www.lancers.jp/work/detail/5412306"
    4  = "https://www.lancers.jp/work/detail/5412467"
    5  = "https://www.lancers.jp/work/detail/5251319"
    6  = "https://www.lancers.jp/work/detail/5412233"
    7  = "https://www.lancers.jp/work/detail/5412194"
    8  = "https://www.lancers.jp/work/detail/5407811"
    9  = "https://www.lancers.jp/work/detail/5412179"
    10 = "https://www.lancers.jp/work/detail/5412487"
    11 = "https://www.lancers.jp/work/detail/5412453"
    12 = "https://www.lancers.jp/work/detail/5412261"
    13 = "https://www.lancers.jp/work/detail/5400683"
    14 = "https://www.lancers.jp/work/detail/5411887"
    15 = "https://www.lancers.jp/work/detail/5412357"
}

$ws.Hyperlinks.Delete()
foreach ($r in 2..15) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $urls[$r])
    # Keep the original "Hyperlink" cell style (blue/underline) rather than
    # whatever ad-hoc style the Add() call might otherwise introduce.
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
